$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-40 down to 16-41
$ws.Rows(15).Insert()

# Populate the new row 15 with the new data record
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44536
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103003
$ws.Cells.Item(15, 10).Value = "Damasco"
$ws.Cells.Item(15, 11).Value = "Castle Brite"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 112
$ws.Cells.Item(15, 14).Value = 13000
$ws.Cells.Item(15, 15).Value = 13000
$ws.Cells.Item(15, 16).Value = 13000
$ws.Cells.Item(15, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(15, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(15, 19).Value = 1300
$ws.Cells.Item(15, 20).Value = 10
